# Refresh the cryptocurrency price (column D) and 1h volume-change
# (column E) snapshot values for rows 2-51 on the active worksheet,
# matching the latest data pull from the source feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "26.791.18" },
    @{ Cell = "E2"; Value = "  -1.75%  " },
    @{ Cell = "D3"; Value = "1.548.63" },
    @{ Cell = "E3"; Value = "  -1.76%  " },
    @{ Cell = "E4"; Value = "  +0.06%  " },
    @{ Cell = "D5"; Value = "204.66" },
    @{ Cell = "E5"; Value = "  -1.62%  " },
    @{ Cell = "D6"; Value = "0.482" },
    @{ Cell = "E6"; Value = "  -1.75%  " },
    @{ Cell = "E7"; Value = "  +0.07%  " },
    @{ Cell = "E8"; Value = "  -1.21%  " },
    @{ Cell = "D9"; Value = "21.35" },
    @{ Cell = "E9"; Value = "  -4.20%  " },
    @{ Cell = "E10"; Value = "  -1.82%  " },
    @{ Cell = "E11"; Value = "  -1.22%  " },
    @{ Cell = "D12"; Value = "1.765.84" },
    @{ Cell = "E12"; Value = "  -1.94%  " },
    @{ Cell = "D13"; Value = "1.548.84" },
    @{ Cell = "E13"; Value = "  -1.70%  " },
    @{ Cell = "E14"; Value = "  -2.72%  " },
    @{ Cell = "D15"; Value = "0.510" },
    @{ Cell = "E15"; Value = "  -2.06%  " },
    @{ Cell = "D16"; Value = "26.783.55" },
    @{ Cell = "E16"; Value = "  -1.84%  " },
    @{ Cell = "D17"; Value = "60.97" },
    @{ Cell = "E17"; Value = "  -2.68%  " },
    @{ Cell = "D18"; Value = "213.62" },
    @{ Cell = "E18"; Value = "  -1.06%  " },
    @{ Cell = "E19"; Value = "  -0.78%  " },
    @{ Cell = "D20"; Value = "0.0₃0682" },
    @{ Cell = "E20"; Value = "  -1.02%  " },
    @{ Cell = "E21"; Value = "  +0.20%  " },
    @{ Cell = "D22"; Value = "4.08" },
    @{ Cell = "E22"; Value = "  -1.53%  " },
    @{ Cell = "D23"; Value = "9.04" },
    @{ Cell = "E23"; Value = "  -4.11%  " },
    @{ Cell = "E24"; Value = "  -0.37%  " },
    @{ Cell = "D25"; Value = "152.65" },
    @{ Cell = "E25"; Value = "  +0.85%  " },
    @{ Cell = "D26"; Value = "6.50" },
    @{ Cell = "E26"; Value = "  -2.89%  " },
    @{ Cell = "D27"; Value = "14.81" },
    @{ Cell = "E27"; Value = "  -1.07%  " },
    @{ Cell = "E28"; Value = "  +0.05%  " },
    @{ Cell = "E29"; Value = "  -2.52%  " },
    @{ Cell = "E30"; Value = "  -0.64%  " },
    @{ Cell = "E31"; Value = "  -3.41%  " },
    @{ Cell = "E32"; Value = "  -0.68%  " },
    @{ Cell = "D33"; Value = "1.356.79" },
    @{ Cell = "E33"; Value = "  -3.43%  " },
    @{ Cell = "E34"; Value = "  -0.90%  " },
    @{ Cell = "D35"; Value = "1.50" },
    @{ Cell = "E35"; Value = "  -4.36%  " },
    @{ Cell = "E36"; Value = "  -0.66%  " },
    @{ Cell = "D37"; Value = "0.912" },
    @{ Cell = "E37"; Value = "  -3.20%  " },
    @{ Cell = "E38"; Value = "  -2.35%  " },
    @{ Cell = "E39"; Value = "  +0.27%  " },
    @{ Cell = "D40"; Value = "0.801" },
    @{ Cell = "E40"; Value = "  -2.37%  " },
    @{ Cell = "E41"; Value = "  +0.07%  " },
    @{ Cell = "E42"; Value = "  -1.23%  " },
    @{ Cell = "D43"; Value = "5.52" },
    @{ Cell = "E43"; Value = "  +3.13%  " },
    @{ Cell = "D44"; Value = "2.19" },
    @{ Cell = "E44"; Value = "  -0.05%  " },
    @{ Cell = "D45"; Value = "1.77" },
    @{ Cell = "E45"; Value = "  -2.34%  " },
    @{ Cell = "D46"; Value = "62.86" },
    @{ Cell = "E46"; Value = "  -1.68%  " },
    @{ Cell = "D47"; Value = "2.31" },
    @{ Cell = "E47"; Value = "  -2.56%  " },
    @{ Cell = "D48"; Value = "1.680.74" },
    @{ Cell = "E48"; Value = "  -1.86%  " },
    @{ Cell = "D49"; Value = "85.78" },
    @{ Cell = "E49"; Value = "  -0.55%  " },
    @{ Cell = "D50"; Value = "0.0508" },
    @{ Cell = "E50"; Value = "  +2.65%  " },
    @{ Cell = "D51"; Value = "0.0₇0962" },
    @{ Cell = "E51"; Value = "  -2.88%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force text storage so purely-numeric-looking values (e.g. "0.482",
    # "204.66") aren't silently reinterpreted as numbers by Excel, then
    # drop back to the default style so no formatting change is introduced.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
